$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (row 11) - local storage / proceso row
# Written in this order so new shared-string entries land in the same
# sequence as the saved workbook (C, then B, then A).
$ws.Range("C11").Value = "juzpassa"
$ws.Range("B11").Value = "jejeje"
$ws.Range("A11").Value = "pepe"
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 12

# Update the active selection to match the saved workbook state
$ws.Range("B18").Select()
